# Natmi following Dr Hou advice
# Recomputed C3 -> Nrp1 ligand/receptor statistics across the ECs/FAPs/sCs
# cluster pairs; the table now has a full 3x3 matrix of sending/target
# clusters (rows 2-10) instead of the previous partial 2x3 matrix (rows 2-7).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "C3"
$ws.Range("C2").Value = "Nrp1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 34.739995
$ws.Range("H2").Value = 104.219985
$ws.Range("I2").Value = 0.1827267341390226
$ws.Range("J2").Value = 0.1827267341390226
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 102.8289443333334
$ws.Range("N2").Value = 308.486833
$ws.Range("O2").Value = 0.5559120396302444
$ws.Range("P2").Value = 0.5559120396302443
$ws.Range("Q2").Value = 3572.277011995279
$ws.Range("R2").Value = 32150.49310795751
$ws.Range("S2").Value = 0.1015799914701975
$ws.Range("T2").Value = 0.1015799914701975

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "C3"
$ws.Range("C3").Value = "Nrp1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 34.739995
$ws.Range("H3").Value = 104.219985
$ws.Range("I3").Value = 0.1827267341390226
$ws.Range("J3").Value = 0.1827267341390226
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 63.66262833333334
$ws.Range("N3").Value = 190.987885
$ws.Range("O3").Value = 0.3441717873742006
$ws.Range("P3").Value = 0.3441717873742006
$ws.Range("Q3").Value = 2211.639389986859
$ws.Range("R3").Value = 19904.75450988172
$ws.Range("S3").Value = 0.06288938668967778
$ws.Range("T3").Value = 0.06288938668967778

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "C3"
$ws.Range("C4").Value = "Nrp1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 34.739995
$ws.Range("H4").Value = 104.219985
$ws.Range("I4").Value = 0.1827267341390226
$ws.Range("J4").Value = 0.1827267341390226
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 18.481835
$ws.Range("N4").Value = 55.445505
$ws.Range("O4").Value = 0.09991617299555507
$ws.Range("P4").Value = 0.09991617299555505
$ws.Range("Q4").Value = 642.0588554908251
$ws.Range("R4").Value = 5778.529699417424
$ws.Range("S4").Value = 0.01825735597914738
$ws.Range("T4").Value = 0.01825735597914738

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "C3"
$ws.Range("C5").Value = "Nrp1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 154.8642143333334
$ws.Range("H5").Value = 464.5926430000001
$ws.Range("I5").Value = 0.8145606273154508
$ws.Range("J5").Value = 0.8145606273154508
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 102.8289443333334
$ws.Range("N5").Value = 308.486833
$ws.Range("O5").Value = 0.5559120396302444
$ws.Range("P5").Value = 0.5559120396302443
$ws.Range("Q5").Value = 15924.52367490774
$ws.Range("R5").Value = 143320.7130741697
$ws.Range("S5").Value = 0.4528240597334236
$ws.Range("T5").Value = 0.4528240597334235

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "C3"
$ws.Range("C6").Value = "Nrp1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 154.8642143333334
$ws.Range("H6").Value = 464.5926430000001
$ws.Range("I6").Value = 0.8145606273154508
$ws.Range("J6").Value = 0.8145606273154508
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 63.66262833333334
$ws.Range("N6").Value = 190.987885
$ws.Range("O6").Value = 0.3441717873742006
$ws.Range("P6").Value = 0.3441717873742006
$ws.Range("Q6").Value = 9859.062919236676
$ws.Range("R6").Value = 88731.56627313007
$ws.Range("S6").Value = 0.2803487870278088
$ws.Range("T6").Value = 0.2803487870278088

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "C3"
$ws.Range("C7").Value = "Nrp1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 154.8642143333334
$ws.Range("H7").Value = 464.5926430000001
$ws.Range("I7").Value = 0.8145606273154508
$ws.Range("J7").Value = 0.8145606273154508
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 18.481835
$ws.Range("N7").Value = 55.445505
$ws.Range("O7").Value = 0.09991617299555507
$ws.Range("P7").Value = 0.09991617299555505
$ws.Range("Q7").Value = 2862.174856713302
$ws.Range("R7").Value = 25759.57371041972
$ws.Range("S7").Value = 0.08138778055421844
$ws.Range("T7").Value = 0.08138778055421843

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "C3"
$ws.Range("C8").Value = "Nrp1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.5157266666666667
$ws.Range("H8").Value = 1.54718
$ws.Range("I8").Value = 0.002712638545526686
$ws.Range("J8").Value = 0.002712638545526686
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 102.8289443333334
$ws.Range("N8").Value = 308.486833
$ws.Range("O8").Value = 0.5559120396302444
$ws.Range("P8").Value = 0.5559120396302443
$ws.Range("Q8").Value = 53.03162869788223
$ws.Range("R8").Value = 477.2846582809401
$ws.Range("S8").Value = 0.00150798842662336
$ws.Range("T8").Value = 0.001507988426623359

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "C3"
$ws.Range("C9").Value = "Nrp1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.5157266666666667
$ws.Range("H9").Value = 1.54718
$ws.Range("I9").Value = 0.002712638545526686
$ws.Range("J9").Value = 0.002712638545526686
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 63.66262833333334
$ws.Range("N9").Value = 190.987885
$ws.Range("O9").Value = 0.3441717873742006
$ws.Range("P9").Value = 0.3441717873742006
$ws.Range("Q9").Value = 32.83251510158889
$ws.Range("R9").Value = 295.4926359143
$ws.Range("S9").Value = 0.0009336136567140713
$ws.Range("T9").Value = 0.0009336136567140715

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "C3"
$ws.Range("C10").Value = "Nrp1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5157266666666667
$ws.Range("H10").Value = 1.54718
$ws.Range("I10").Value = 0.002712638545526686
$ws.Range("J10").Value = 0.002712638545526686
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 18.481835
$ws.Range("N10").Value = 55.445505
$ws.Range("O10").Value = 0.09991617299555507
$ws.Range("P10").Value = 0.09991617299555505
$ws.Range("Q10").Value = 9.531575158433334
$ws.Range("R10").Value = 85.7841764259
$ws.Range("S10").Value = 0.0002710364621892552
$ws.Range("T10").Value = 0.0002710364621892552

